$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = @(
    @(8, 8),
    @(8, 9),
    @(7, 7),
    @(10, 10),
    @(9, 9),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(7, 7),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(7, 8),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(5, 5)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $data[$r][0]
    $ws.Cells.Item($row, 10).Value = $data[$r][1]
}
